# Apply updated odds values to Sheet1 as described by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.67

# Row 10 updates
$ws.Range("G10").Value = 1.16
$ws.Range("H10").Value = 5.5
$ws.Range("I10").Value = 16
$ws.Range("J10").Value = 1.53
$ws.Range("K10").Value = 2.6
$ws.Range("L10").Value = 11.25
$ws.Range("O10").Value = 1.11
$ws.Range("P10").Value = 4.55
$ws.Range("S10").Value = 1.25
$ws.Range("T10").Value = 3.56
$ws.Range("U10").Value = 2.28
$ws.Range("V10").Value = 1.6
$ws.Range("W10").Value = 6.2
$ws.Range("Y10").Value = 8.25
$ws.Range("Z10").Value = 5.4
$ws.Range("AA10").Value = 9.25
$ws.Range("AB10").Value = 27
$ws.Range("AC10").Value = 14
$ws.Range("AD10").Value = 10.5
$ws.Range("AE10").Value = 24
$ws.Range("AF10").Value = 100
$ws.Range("AH10").Value = 35
$ws.Range("AI10").Value = 120
$ws.Range("AJ10").Value = 40
$ws.Range("AN10").Value = 2.92
$ws.Range("AO10").Value = 4.7
$ws.Range("AQ10").Value = 11
$ws.Range("AT10").Value = 3.35
$ws.Range("AV10").Value = 90
$ws.Range("AW10").Value = 14.5
